$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 1456
$ws.Cells.Item(5, 6).Value = 7102
$ws.Cells.Item(6, 6).Value = 569
$ws.Cells.Item(7, 6).Value = 1060
$ws.Cells.Item(9, 6).Value = 4664
$ws.Cells.Item(10, 6).Value = 6880
$ws.Cells.Item(12, 6).Value = 243
$ws.Cells.Item(13, 6).Value = 1439
$ws.Cells.Item(14, 6).Value = 832
$ws.Cells.Item(15, 6).Value = 136
$ws.Cells.Item(17, 6).Value = 1145
$ws.Cells.Item(19, 6).Value = 142
$ws.Cells.Item(21, 6).Value = 195
$ws.Cells.Item(23, 6).Value = 1098
$ws.Cells.Item(24, 6).Value = 542
$ws.Cells.Item(25, 6).Value = 41
$ws.Cells.Item(26, 6).Value = 1181
$ws.Cells.Item(28, 6).Value = 128
$ws.Cells.Item(31, 6).Value = 133
$ws.Cells.Item(33, 6).Value = 20
$ws.Cells.Item(34, 6).Value = 38
$ws.Cells.Item(36, 6).Value = 29
$ws.Cells.Item(37, 6).Value = 537
$ws.Cells.Item(38, 6).Value = 399
$ws.Cells.Item(39, 6).Value = 59
$ws.Cells.Item(40, 6).Value = 57
$ws.Cells.Item(41, 6).Value = 337
$ws.Cells.Item(43, 6).Value = 550
$ws.Cells.Item(44, 6).Value = 71
$ws.Cells.Item(45, 6).Value = 127
$ws.Cells.Item(47, 6).Value = 10

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 23
$ws.Cells.Item(3, 6).Value = 23
$ws.Cells.Item(13, 6).Value = 23
$ws.Cells.Item(18, 6).Value = 552
$ws.Cells.Item(27, 6).Value = 619
$ws.Cells.Item(32, 6).Value = 830
$ws.Cells.Item(34, 6).Value = 593
$ws.Cells.Item(41, 6).Value = 137
$ws.Cells.Item(44, 6).Value = 68

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 719
$ws.Cells.Item(5, 6).Value = 849
$ws.Cells.Item(6, 6).Value = 637
$ws.Cells.Item(8, 6).Value = 1426
$ws.Cells.Item(9, 6).Value = 2254

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 23
$ws.Cells.Item(3, 6).Value = 719
$ws.Cells.Item(6, 6).Value = 637
$ws.Cells.Item(7, 6).Value = 637
$ws.Cells.Item(11, 6).Value = 7103
$ws.Cells.Item(12, 6).Value = 569
$ws.Cells.Item(14, 6).Value = 4664
$ws.Cells.Item(15, 6).Value = 6880
$ws.Cells.Item(16, 6).Value = 1439
$ws.Cells.Item(17, 6).Value = 23
$ws.Cells.Item(20, 6).Value = 552
$ws.Cells.Item(21, 6).Value = 136
$ws.Cells.Item(22, 6).Value = 1426
$ws.Cells.Item(23, 6).Value = 2254
$ws.Cells.Item(26, 6).Value = 1145
$ws.Cells.Item(27, 6).Value = 142
$ws.Cells.Item(31, 6).Value = 1098
$ws.Cells.Item(32, 6).Value = 619
$ws.Cells.Item(33, 6).Value = 542
$ws.Cells.Item(34, 6).Value = 41
$ws.Cells.Item(35, 6).Value = 1181
$ws.Cells.Item(39, 6).Value = 133
$ws.Cells.Item(41, 6).Value = 830
$ws.Cells.Item(44, 6).Value = 537
$ws.Cells.Item(45, 6).Value = 59
$ws.Cells.Item(48, 6).Value = 71
